# Weekly Fruta/Hortaliza update: insert a new price record for
# "Feria Lagunitas de Puerto Montt" - Piña (Caramelo, Segunda) above the
# existing row 473, shifting all subsequent records down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 473 (existing rows 473:504 shift down to 474:505)
$ws.Rows.Item(473).Insert()

# Populate the new row with the latest observation
$ws.Cells.Item(473, 1).Value  = 4
$ws.Cells.Item(473, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(473, 3).Value  = "Los Lagos"
$ws.Cells.Item(473, 4).Value  = 45265
$ws.Cells.Item(473, 5).Value  = 10
$ws.Cells.Item(473, 6).Value  = "Fruta"
$ws.Cells.Item(473, 7).Value  = 100108
$ws.Cells.Item(473, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(473, 9).Value  = 100108005
$ws.Cells.Item(473, 10).Value = "Piña"
$ws.Cells.Item(473, 11).Value = "Caramelo"
$ws.Cells.Item(473, 12).Value = "Segunda"
$ws.Cells.Item(473, 13).Value = 120
$ws.Cells.Item(473, 14).Value = 28000
$ws.Cells.Item(473, 15).Value = 28000
$ws.Cells.Item(473, 16).Value = 28000
$ws.Cells.Item(473, 17).Value = "$/caja 14 unidades"
$ws.Cells.Item(473, 18).Value = "Ecuador"
$ws.Cells.Item(473, 19).Value = 2000
$ws.Cells.Item(473, 20).Value = 14
